$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Remove slides with SlideId 260, 261 and 257 (the three trailing slides
#    that were deleted from the deck: two empty stub slides and the big
#    "insights" content slide). Slides 258 and 259 remain.
# ---------------------------------------------------------------------------
foreach ($sid in 260, 261, 257) {
    $slide = $p.Slides.FindBySlideID($sid)
    if ($slide -ne $null) {
        $slide.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2) Refresh the cached "datetimeFigureOut" footer date from 20/08/2024 to
#    26/08/2024 everywhere it is shown: the Slide Master and every Custom
#    Layout's Date placeholder.
# ---------------------------------------------------------------------------
$oldDate = "20/08/2024"
$newDate = "26/08/2024"

function Update-DateShape($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}
